# Add two new columns, I ("I0") and J ("IF"), to the stats table on the
# active sheet. I is a constant 1 for every data row; J duplicates the
# existing H ("IP") value for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in I1/J1, matching H1's style ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-40): J = copy of existing H values, I = constant 1 ---
$ws.Range("H2:H40").Copy()
$ws.Range("J2:J40").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("I2:I40").Value = 1

$excel.CutCopyMode = 0
